# Apply cached-value updates (e.g. refreshed market-board prices) to the
# per-job profit tables, matching the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Cells.Item(80, 8).Value = 903.4
$ws.Cells.Item(80, 9).Value = 518.7778
$ws.Cells.Item(80, 10).Value = 1218.091
$ws.Cells.Item(80, 11).Value = 1556.3334
$ws.Cells.Item(80, 12).Value = 3654.273
$ws.Cells.Item(80, 13).Value = -558.3334
$ws.Cells.Item(80, 14).Value = -5650.272999999999
# Row 83
$ws.Cells.Item(83, 8).Value = 903.4
$ws.Cells.Item(83, 9).Value = 518.7778
$ws.Cells.Item(83, 10).Value = 1218.091
$ws.Cells.Item(83, 11).Value = 4669.000199999999
$ws.Cells.Item(83, 12).Value = 10962.819
$ws.Cells.Item(83, 13).Value = 322.9998000000005
$ws.Cells.Item(83, 14).Value = -20946.819
# Row 100
$ws.Cells.Item(100, 8).Value = 32626.438
$ws.Cells.Item(100, 9).Value = 51189.3
$ws.Cells.Item(100, 10).Value = 1688.3334
$ws.Cells.Item(100, 11).Value = 51189.3
$ws.Cells.Item(100, 12).Value = 1688.3334
$ws.Cells.Item(100, 13).Value = -50648.3
$ws.Cells.Item(100, 14).Value = -2770.3334
# Row 107
$ws.Cells.Item(107, 8).Value = 417.16666
$ws.Cells.Item(107, 9).Value = 430.52942
$ws.Cells.Item(107, 11).Value = 430.52942
$ws.Cells.Item(107, 13).Value = 1489.47058
# Row 141
$ws.Cells.Item(141, 8).Value = 5758.8335
$ws.Cells.Item(141, 9).Value = 2245.2632
$ws.Cells.Item(141, 10).Value = 11827.728
$ws.Cells.Item(141, 11).Value = 6735.7896
$ws.Cells.Item(141, 12).Value = 35483.18399999999
$ws.Cells.Item(141, 13).Value = -1555.7896
$ws.Cells.Item(141, 14).Value = -45843.18399999999

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Cells.Item(22, 8).Value = 4102
$ws.Cells.Item(22, 9).Value = 1040.75
$ws.Cells.Item(22, 10).Value = 9000
$ws.Cells.Item(22, 11).Value = 1040.75
$ws.Cells.Item(22, 12).Value = 9000
$ws.Cells.Item(22, 13).Value = -741.75
$ws.Cells.Item(22, 14).Value = -9598
# Row 45
$ws.Cells.Item(45, 8).Value = 1229.6
$ws.Cells.Item(45, 9).Value = 1162
$ws.Cells.Item(45, 10).Value = 1500
$ws.Cells.Item(45, 11).Value = 1162
$ws.Cells.Item(45, 12).Value = 1500
$ws.Cells.Item(45, 13).Value = -785
$ws.Cells.Item(45, 14).Value = -2254
# Row 105
$ws.Cells.Item(105, 8).Value = 28000
$ws.Cells.Item(105, 10).Value = 28000
$ws.Cells.Item(105, 12).Value = 28000
$ws.Cells.Item(105, 14).Value = -34988
# Row 132
$ws.Cells.Item(132, 8).Value = 7378.2705
$ws.Cells.Item(132, 9).Value = 5514.88
$ws.Cells.Item(132, 10).Value = 11260.333
$ws.Cells.Item(132, 11).Value = 16544.64
$ws.Cells.Item(132, 12).Value = 33780.999
$ws.Cells.Item(132, 13).Value = -14014.64
$ws.Cells.Item(132, 14).Value = -38840.999

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Cells.Item(60, 8).Value = 48000
$ws.Cells.Item(60, 10).Value = 48000
$ws.Cells.Item(60, 12).Value = 48000
$ws.Cells.Item(60, 14).Value = -49198
# Row 86
$ws.Cells.Item(86, 8).Value = 273026.1
$ws.Cells.Item(86, 9).Value = 2055.1304
$ws.Cells.Item(86, 10).Value = 718192.7
$ws.Cells.Item(86, 11).Value = 2055.1304
$ws.Cells.Item(86, 12).Value = 718192.7
$ws.Cells.Item(86, 13).Value = -932.1304
$ws.Cells.Item(86, 14).Value = -720438.7
# Row 89
$ws.Cells.Item(89, 8).Value = 273026.1
$ws.Cells.Item(89, 9).Value = 2055.1304
$ws.Cells.Item(89, 10).Value = 718192.7
$ws.Cells.Item(89, 11).Value = 10275.652
$ws.Cells.Item(89, 12).Value = 3590963.5
$ws.Cells.Item(89, 13).Value = -4659.652
$ws.Cells.Item(89, 14).Value = -3602195.5
# Row 94
$ws.Cells.Item(94, 8).Value = 678.3
$ws.Cells.Item(94, 9).Value = 457
$ws.Cells.Item(94, 10).Value = 1089.2858
$ws.Cells.Item(94, 11).Value = 457
$ws.Cells.Item(94, 12).Value = 1089.2858
$ws.Cells.Item(94, 13).Value = -6
$ws.Cells.Item(94, 14).Value = -1991.2858
# Row 105
$ws.Cells.Item(105, 8).Value = 2950.0334
$ws.Cells.Item(105, 9).Value = 3044.6155
$ws.Cells.Item(105, 10).Value = 2877.7058
$ws.Cells.Item(105, 11).Value = 3044.6155
$ws.Cells.Item(105, 12).Value = 2877.7058
$ws.Cells.Item(105, 13).Value = -1297.6155
$ws.Cells.Item(105, 14).Value = -6371.7058

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Cells.Item(5, 8).Value = 274.2143
$ws.Cells.Item(5, 9).Value = 149.33333
$ws.Cells.Item(5, 10).Value = 367.875
$ws.Cells.Item(5, 11).Value = 149.33333
$ws.Cells.Item(5, 12).Value = 367.875
$ws.Cells.Item(5, 13).Value = -37.33332999999999
$ws.Cells.Item(5, 14).Value = -591.875
# Row 122
$ws.Cells.Item(122, 8).Value = 1479.1428
$ws.Cells.Item(122, 9).Value = 1632.5
$ws.Cells.Item(122, 10).Value = 1274.6666
$ws.Cells.Item(122, 11).Value = 4897.5
$ws.Cells.Item(122, 12).Value = 3823.9998
$ws.Cells.Item(122, 13).Value = -2447.5
$ws.Cells.Item(122, 14).Value = -8723.9998
# Row 133
$ws.Cells.Item(133, 8).Value = 42469
$ws.Cells.Item(133, 9).Value = 12345
$ws.Cells.Item(133, 10).Value = 50000
$ws.Cells.Item(133, 11).Value = 12345
$ws.Cells.Item(133, 12).Value = 50000
$ws.Cells.Item(133, 13).Value = -9815
$ws.Cells.Item(133, 14).Value = -55060

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 54.846153
$ws.Cells.Item(8, 9).Value = 54.846153
$ws.Cells.Item(8, 11).Value = 164.538459
$ws.Cells.Item(8, 13).Value = -25.53845899999999
# Row 80
$ws.Cells.Item(80, 8).Value = 3119.6667
$ws.Cells.Item(80, 9).Value = 2088
$ws.Cells.Item(80, 11).Value = 6264
$ws.Cells.Item(80, 13).Value = -5328
# Row 83
$ws.Cells.Item(83, 8).Value = 3119.6667
$ws.Cells.Item(83, 9).Value = 2088
$ws.Cells.Item(83, 11).Value = 18792
$ws.Cells.Item(83, 13).Value = -14112
# Row 123
$ws.Cells.Item(123, 8).Value = 1000
$ws.Cells.Item(123, 9).Value = 1000
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 3000
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).Value = -550
$ws.Cells.Item(123, 14).ClearContents()  # N123
# Row 131
$ws.Cells.Item(131, 8).Value = 631.98
$ws.Cells.Item(131, 9).Value = 277.13635
$ws.Cells.Item(131, 10).Value = 910.7857
$ws.Cells.Item(131, 11).Value = 831.40905
$ws.Cells.Item(131, 12).Value = 2732.3571
$ws.Cells.Item(131, 13).Value = 4208.59095
$ws.Cells.Item(131, 14).Value = -12812.3571
# Row 132
$ws.Cells.Item(132, 8).Value = 1059.5652
$ws.Cells.Item(132, 9).Value = 897.3333
$ws.Cells.Item(132, 10).Value = 1363.75
$ws.Cells.Item(132, 11).Value = 8075.9997
$ws.Cells.Item(132, 12).Value = 12273.75
$ws.Cells.Item(132, 13).Value = -5545.9997
$ws.Cells.Item(132, 14).Value = -17333.75

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Cells.Item(13, 8).Value = 412.5
$ws.Cells.Item(13, 9).Value = 216.66667
$ws.Cells.Item(13, 11).Value = 216.66667
$ws.Cells.Item(13, 13).Value = -77.66667000000001
# Row 80
$ws.Cells.Item(80, 8).Value = 3235.5557
$ws.Cells.Item(80, 9).Value = 2103.3333
$ws.Cells.Item(80, 10).Value = 5500
$ws.Cells.Item(80, 11).Value = 2103.3333
$ws.Cells.Item(80, 12).Value = 5500
$ws.Cells.Item(80, 13).Value = -1105.3333
$ws.Cells.Item(80, 14).Value = -7496
# Row 83
$ws.Cells.Item(83, 8).Value = 3235.5557
$ws.Cells.Item(83, 9).Value = 2103.3333
$ws.Cells.Item(83, 10).Value = 5500
$ws.Cells.Item(83, 11).Value = 10516.6665
$ws.Cells.Item(83, 12).Value = 27500
$ws.Cells.Item(83, 13).Value = -5524.666499999999
$ws.Cells.Item(83, 14).Value = -37484

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 111
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()  # N111
